$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete two rows to shrink the table from 18 to 16 data rows ---
# Delete higher row index first to keep lower indices stable
$ws.Rows(15).Delete()
$ws.Rows(11).Delete()

# --- Step 2: update the title cell ---
$ws.Range("B1").Value2 = "Partner Management (MISP and E-KYC/Auth Partners)"

# --- Step 3: rewrite all data rows (B3:F18) with the updated requirements content ---
# Row 3: Sr No. 1
$ws.Range("B3").Value2 = 1
$ws.Range("C3").Value2 = "Kernel"
$ws.Range("D3").Value2 = "MISP ID Generation"
$ws.Range("E3").Value2 = "1. Generate MISP ID as per below logic`na. MISP ID should be of 3 digits (Configurable)`nb. MISP ID should be generated sequentially`nc. MISP ID should be generated incrementally for every request"
$ws.Range("F3").Value2 = "Component already exist as TSP ID generator"
$ws.Rows(3).RowHeight = 58

# Row 4: Sr No. 2
$ws.Range("B4").Value2 = 2
$ws.Range("C4").Value2 = "Kernel"
$ws.Range("D4").Value2 = "MISP License Key Generation"
$ws.Range("E4").Value2 = "1. Generate a License Key as per below logic`na. License Key generation to follow random pattern`nb. License Key should be alphanumeric`nc. Length should be 8 digits (Configurable)`nd. Should be mapped to an expiry"
$ws.Range("F4").ClearContents()
$ws.Rows(4).RowHeight = 72.5

# Row 5: Sr No. 3
$ws.Range("B5").Value2 = 3
$ws.Range("C5").Value2 = "Kernel"
$ws.Range("D5").Value2 = "MISP License Key Pattern Validation"
$ws.Range("E5").Value2 = "1. Validate length of a License Key as configured and respond as mentioned below`na. If found valid, respond with ""VALID""`nb. if found invalid, respond with ""INVALID"""
$ws.Range("F5").ClearContents()
$ws.Rows(5).RowHeight = 58

# Row 6: Sr No. 4
$ws.Range("B6").Value2 = 4
$ws.Range("C6").Value2 = "Admin"
$ws.Range("D6").Value2 = "MSIP License Key Expiry Validation"
$ws.Range("E6").Value2 = "1. Validate status of Lisence Key and respond as mentioned below`na. If found expired, respond with ""Your License Key is EXPIRED. Please regenrate a new License Key""`nb. If found temporarily sespended, respond with ""Your License Key is temporarily SUSPENDED. Please contact MOSIP Administration""`nc. If found permanently blocked, respond with ""Your License Key is BLOCKED. Please contact MOSIP Administration"""
$ws.Range("F6").ClearContents()
$ws.Rows(6).RowHeight = 101.5

# Row 7: Sr No. 5
$ws.Range("B7").Value2 = 5
$ws.Range("C7").Value2 = "Admin"
$ws.Range("D7").Value2 = "MISP Registration"
$ws.Range("E7").Value2 = "1. Receive request to register a MISP with follwing parameters`na. MISP Name`nb. MISP Contact Name`nc. MISP Phone`nd. MISP Email ID`n2. Issue and Map MISP ID`n3. Issue and Map Lisence Key`n4. Store the MISP in MOSIP"
$ws.Range("F7").ClearContents()
$ws.Rows(7).RowHeight = 116

# Row 8: Sr No. 6
$ws.Range("B8").Value2 = 6
$ws.Range("C8").Value2 = "Kernel"
$ws.Range("D8").Value2 = "Partner ID Generation"
$ws.Range("E8").Value2 = "1. Generate Partner ID as per below logic`na. Partner ID should be of 4 digits (Configurable)`nb. Partner ID should be generated sequentially`nc. Partner ID should be generated incrementally for every request"
$ws.Range("F8").ClearContents()
$ws.Rows(8).RowHeight = 58

# Row 9: Sr No. 7
$ws.Range("B9").Value2 = 7
$ws.Range("C9").Value2 = "Kernel"
$ws.Range("D9").Value2 = "Partner ID Validation"
$ws.Range("E9").Value2 = "1. Validate length of a Partner ID as configured and respond as mentioned below`na. If found valid, respond with ""VALID""`nb. if found invalid, respond with ""INVALID"""
$ws.Range("F9").ClearContents()
$ws.Rows(9).RowHeight = 58

# Row 10: Sr No. 8
$ws.Range("B10").Value2 = 8
$ws.Range("C10").Value2 = "Kernel "
$ws.Range("D10").Value2 = "Map Policies to Partners"
$ws.Range("E10").Value2 = "1. Map following Policies to Partners`na. Auth Policies ( can be Mandatory/Non-Mandatory)`n     1. OTP Trigger `n     2. OTP Authentication`n     3. Demo Authentication `n     4. Biometric Authentication - FMR Data Match `n     5. Biometric Authentication - IIR Data Match  `n     6. Biometric Authentication - FID Data Match `nb. E-Kyc Policies (can be Required/Not Required)`n    1. eKYC - all combinations of eKYC demo fields "
$ws.Range("F10").ClearContents()
$ws.Rows(10).RowHeight = 145

# Row 11: Sr No. 10
$ws.Range("B11").Value2 = 10
$ws.Range("C11").Value2 = "Admin"
$ws.Range("D11").Value2 = "Retrieve Policies based on Partner ID"
$ws.Range("E11").Value2 = "1. Receive request to retreive policies based on Partner ID`n2. Respond appropirately if Partner ID does not exist"
$ws.Range("F11").ClearContents()
$ws.Rows(11).RowHeight = 29

# Row 12: Sr No. 11
$ws.Range("B12").Value2 = 11
$ws.Range("C12").Value2 = "Admin"
$ws.Range("D12").Value2 = "Partner Registration"
$ws.Range("E12").Value2 = "1. Receive request to register a Partner with follwing parameters`na. Partner Name`nb. Partner Contact Name`nc. Partner Phone`nd. Partner Email ID`n2. Issue and Map Partner ID`n3. Map Policies to the Partner`na. Multiple Policies can be mapped to a Partner`nb. A Partner can have a policy for both Auth and E-KYC`n4. Store the Partner in MOSIP"
$ws.Range("F12").ClearContents()
$ws.Rows(12).RowHeight = 145

# Row 13: Sr No. 12
$ws.Range("B13").Value2 = 12
$ws.Range("C13").Value2 = "Admin"
$ws.Range("D13").Value2 = "MISP - Partner Mapping"
$ws.Range("E13").Value2 = "1. Receive a request to map MISP to a Partner with MISP ID and Partner ID as Input`n2. There can be a many-to-mapping between MISPs and Partners"
$ws.Range("F13").ClearContents()
$ws.Rows(13).RowHeight = 43.5

# Row 14: Sr No. 13
$ws.Range("B14").Value2 = 13
$ws.Range("C14").Value2 = "Admin"
$ws.Range("D14").Value2 = "Partner Certificate Validation"
$ws.Range("E14").Value2 = "1. Upload Digital Certificate on Admin Portal for a Partner`n2. Verify CA Authority of the certificate`n3. Sign the certificate with MOSIP Certificate`n4. Respond to the source with the re-issued certificate`n5. Certificate will be uploaded by the MOSIP admin. The Registered Partner will send the certificate to the MOSIP Admin through ofline process. Re-issued certificate will be sent to the Partner by MOSIP admin through notification/offline process`n6. Private key to change priodically as per the Key Rotation Policy set by admin"
$ws.Range("F14").ClearContents()
$ws.Rows(14).RowHeight = 145

# Row 15: Sr No. 15
$ws.Range("B15").Value2 = 15
$ws.Range("C15").Value2 = "Admin"
$ws.Range("D15").Value2 = "Distribution of Public Key to Partners"
$ws.Range("E15").Value2 = "1. Distribute Public Key to Partners for encrypting the Auth Request befoe sending it to the MOSIP`n2. Public key needs to be distributed priodically whenever the correspinding Private Key is rotated"
$ws.Range("F15").ClearContents()
$ws.Rows(15).RowHeight = 58

# Row 16: Sr No. 16
$ws.Range("B16").Value2 = 16
$ws.Range("C16").Value2 = "Admin"
$ws.Range("D16").Value2 = "Device Registration"
$ws.Range("E16").Value2 = "TBD"
$ws.Range("F16").Value2 = "Yet to analyzed"

# Row 17: Sr No. 17
$ws.Range("B17").Value2 = 17
$ws.Range("C17").Value2 = "Admin"
$ws.Range("D17").Value2 = "Device Provider Registration"
$ws.Range("E17").Value2 = "TBD"
$ws.Range("F17").Value2 = "Yet to analyzed"

# Row 18: Sr No. 18
$ws.Range("B18").Value2 = 18
$ws.Range("C18").Value2 = "Admin"
$ws.Range("D18").Value2 = "RD Service Registration"
$ws.Range("E18").Value2 = "TBD"
$ws.Range("F18").Value2 = "Yet to analyzed"

# --- Step 4: update view / selection state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("E15").Select()

Write-Output "done"